$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'69.739.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.08%  '
$ws.Range("D3").Formula = "'3.383.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.01%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Formula = "'190.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.13%  '
$ws.Range("D6").Formula = "'593.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.15%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Formula = "'0.610"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Formula = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +1.94%  '
$ws.Range("D10").Formula = "'6.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.88%  '
$ws.Range("D11").Formula = "'0.418"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.60%  '
$ws.Range("D12").Formula = "'3.972.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.19%  '
$ws.Range("E13").Value = '  -0.73%  '
$ws.Range("D14").Formula = "'28.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.90%  '
$ws.Range("D15").Formula = "'69.689.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.00%  '
$ws.Range("E16").Value = '  +1.33%  '
$ws.Range("D17").Formula = "'3.398.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.53%  '
$ws.Range("D18").Formula = "'451.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +14.61%  '
$ws.Range("D19").Formula = "'5.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("D20").Formula = "'13.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.60%  '
$ws.Range("D21").Formula = "'7.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.63%  '
$ws.Range("D22").Formula = "'76.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.47%  '
$ws.Range("D23").Formula = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Formula = "'0.521"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.81%  '
$ws.Range("E25").Value = '  +3.20%  '
$ws.Range("E26").Value = '  +2.20%  '
$ws.Range("D27").Formula = "'9.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.46%  '
$ws.Range("D28").Formula = "'0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("E29").Value = '  +3.12%  '
$ws.Range("D30").Formula = "'23.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.49%  '
$ws.Range("D31").Formula = "'5.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.41%  '
$ws.Range("E32").Value = '  +2.17%  '
$ws.Range("D33").Formula = "'6.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.28%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").Formula = "'1.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.58%  '
$ws.Range("D36").Formula = "'164.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.59%  '
$ws.Range("E37").Value = '  +2.51%  '
$ws.Range("D38").Formula = "'27.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.39%  '
$ws.Range("E39").Value = '  +0.90%  '
$ws.Range("D40").Formula = "'4.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.42%  '
$ws.Range("D41").Formula = "'6.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.91%  '
$ws.Range("D42").Formula = "'2.748.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.10%  '
$ws.Range("D43").Formula = "'2.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.35%  '
$ws.Range("D44").Formula = "'25.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.31%  '
$ws.Range("D45").Formula = "'0.0688"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("D46").Formula = "'41.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.14%  '
$ws.Range("D47").Formula = "'339.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.54%  '
$ws.Range("E48").Value = '  +2.41%  '
$ws.Range("D49").Formula = "'32.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.25%  '
$ws.Range("E50").Value = '  +4.70%  '
$ws.Range("D51").Formula = "'6.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.55%  '
